$d = $word.ActiveDocument

function Assert-ParagraphText($Paragraph, $Expected) {
    # Range.Text includes the trailing paragraph mark (CR); strip only that.
    $actual = $Paragraph.Range.Text.TrimEnd("`r")
    if ($actual -ne $Expected) {
        throw "Unexpected paragraph text: expected [$Expected] got [$actual]"
    }
}

function Set-ParagraphXml($Paragraph, $InnerXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        $InnerXml +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    [void]$Paragraph.Range.InsertXML($pkg)
}

# --- Paragraph 1: "{{ indicateur }}" -> "{{ indicateur_top }}" split over 3 runs,
#     restyled from Titre3/numbered heading to a centered Corpsdetexte "link" style.
$runRpr = '<w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'
$p1xml = '<w:p>' +
    '<w:pPr>' +
        '<w:pStyle w:val="Corpsdetexte"/>' +
        '<w:spacing w:lineRule="auto" w:line="120"/>' +
        '<w:jc w:val="center"/>' +
        '<w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r>' + $runRpr + '<w:t>{{ indicateur_</w:t></w:r>' +
    '<w:r>' + $runRpr + '<w:t>top</w:t></w:r>' +
    '<w:r>' + $runRpr + '<w:t xml:space="preserve"> }}</w:t></w:r>' +
    '</w:p>'
$para1 = $d.Paragraphs.Item(1)
Assert-ParagraphText $para1 "{{ indicateur }}"
Set-ParagraphXml $para1 $p1xml

# --- Paragraph 2: "pour la mesure " -> same restyle, single run.
$p2xml = '<w:p>' +
    '<w:pPr>' +
        '<w:pStyle w:val="Corpsdetexte"/>' +
        '<w:spacing w:lineRule="auto" w:line="120"/>' +
        '<w:jc w:val="center"/>' +
        '<w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r>' + $runRpr + '<w:t xml:space="preserve">pour la mesure </w:t></w:r>' +
    '</w:p>'
$para2 = $d.Paragraphs.Item(2)
Assert-ParagraphText $para2 "pour la mesure "
Set-ParagraphXml $para2 $p2xml

# --- Paragraph 3: "{{ mesure }}" keeps Titre3 + numbering; only the spacing gains
#     w:line="120" (same "auto" rule). Use the ParagraphFormat API directly instead of
#     InsertXML so the rest of pPr/rPr is left completely untouched.
$p3 = $d.Paragraphs.Item(3)
Assert-ParagraphText $p3 "{{ mesure }}"
$p3.Range.ParagraphFormat.LineSpacingRule = 5
$p3.Range.ParagraphFormat.LineSpacing = 6

Write-Host "done"
